$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, pushing the existing rows 118-184 down to 119-185
$ws.Rows("118").Insert()

# Populate the newly inserted row 118 with the new weekly record
$ws.Range("A118").Value = 4
$ws.Range("B118").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C118").Value = "Los Lagos"
$ws.Range("D118").Value = 44572
$ws.Range("E118").Value = 10
$ws.Range("F118").Value = 100112017
$ws.Range("G118").Value = "Apio"
$ws.Range("H118").Value = "Americana (o)"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 35
$ws.Range("K118").Value = 12000
$ws.Range("L118").Value = 12000
$ws.Range("M118").Value = 12000
$ws.Range("N118").Value = "$/docena de matas"
$ws.Range("O118").Value = "Región de Coquimbo"
$ws.Range("P118").Value = 2000
$ws.Range("Q118").Value = 6
$ws.Range("R118").Value = "Hortaliza"
